$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 7 (columns D, J, K, L, M, P)
$r2 = @{ D = $ws.Range("D2").Value2; J = $ws.Range("J2").Value2; K = $ws.Range("K2").Value2; L = $ws.Range("L2").Value2; M = $ws.Range("M2").Value2; P = $ws.Range("P2").Value2 }
$r7 = @{ D = $ws.Range("D7").Value2; J = $ws.Range("J7").Value2; K = $ws.Range("K7").Value2; L = $ws.Range("L7").Value2; M = $ws.Range("M7").Value2; P = $ws.Range("P7").Value2 }

$ws.Range("D2").Value2 = $r7.D
$ws.Range("J2").Value2 = $r7.J
$ws.Range("K2").Value2 = $r7.K
$ws.Range("L2").Value2 = $r7.L
$ws.Range("M2").Value2 = $r7.M
$ws.Range("P2").Value2 = $r7.P

$ws.Range("D7").Value2 = $r2.D
$ws.Range("J7").Value2 = $r2.J
$ws.Range("K7").Value2 = $r2.K
$ws.Range("L7").Value2 = $r2.L
$ws.Range("M7").Value2 = $r2.M
$ws.Range("P7").Value2 = $r2.P

# Swap values between row 3 and row 6 (columns D, J, K, L, M, P)
$r3 = @{ D = $ws.Range("D3").Value2; J = $ws.Range("J3").Value2; K = $ws.Range("K3").Value2; L = $ws.Range("L3").Value2; M = $ws.Range("M3").Value2; P = $ws.Range("P3").Value2 }
$r6 = @{ D = $ws.Range("D6").Value2; J = $ws.Range("J6").Value2; K = $ws.Range("K6").Value2; L = $ws.Range("L6").Value2; M = $ws.Range("M6").Value2; P = $ws.Range("P6").Value2 }

$ws.Range("D3").Value2 = $r6.D
$ws.Range("J3").Value2 = $r6.J
$ws.Range("K3").Value2 = $r6.K
$ws.Range("L3").Value2 = $r6.L
$ws.Range("M3").Value2 = $r6.M
$ws.Range("P3").Value2 = $r6.P

$ws.Range("D6").Value2 = $r3.D
$ws.Range("J6").Value2 = $r3.J
$ws.Range("K6").Value2 = $r3.K
$ws.Range("L6").Value2 = $r3.L
$ws.Range("M6").Value2 = $r3.M
$ws.Range("P6").Value2 = $r3.P

# Swap values between row 4 and row 8 (columns D, J, M, P)
$r4 = @{ D = $ws.Range("D4").Value2; J = $ws.Range("J4").Value2; M = $ws.Range("M4").Value2; P = $ws.Range("P4").Value2 }
$r8 = @{ D = $ws.Range("D8").Value2; J = $ws.Range("J8").Value2; M = $ws.Range("M8").Value2; P = $ws.Range("P8").Value2 }

$ws.Range("D4").Value2 = $r8.D
$ws.Range("J4").Value2 = $r8.J
$ws.Range("M4").Value2 = $r8.M
$ws.Range("P4").Value2 = $r8.P

$ws.Range("D8").Value2 = $r4.D
$ws.Range("J8").Value2 = $r4.J
$ws.Range("M8").Value2 = $r4.M
$ws.Range("P8").Value2 = $r4.P
